$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update row 2's "what" cell: Instructor -> Assistant Instructor
$ws.Range("D2").Value = "Assistant Instructor"

# Insert a new row for the 2022 activity, pushing everything below down by one
$ws.Rows(3).Insert()

$ws.Range("A3").Value = 2022
$ws.Range("B3").Value = "Online"
$ws.Range("D3").Value = "Instructor and content developer"
$ws.Range("E3").Value = "An introduction to R and Statistics for Ecologists"

# Update the Benchling row's "when" value (2020-present -> 2020-2022)
$ws.Range("A7").Value = "2020-2022"

# Update the selection to match the post-edit view
$ws.Range("A9:XFD9").Select()
